$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'39.492.94"
$ws.Range("D3").Value = "'2.165.08"
$ws.Range("E3").Value = "  +2.87%  "
$ws.Range("D5").Value = "'228.68"
$ws.Range("E5").Value = "  +0.03%  "
$ws.Range("D6").Value = "'0.623"
$ws.Range("E6").Value = "  +1.09%  "
$ws.Range("D7").Value = "'63.88"
$ws.Range("E7").Value = "  +2.61%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +1.92%  "
$ws.Range("D10").Value = "'0.0855"
$ws.Range("E10").Value = "  +1.39%  "
$ws.Range("E11").Value = "  -0.04%  "
$ws.Range("D12").Value = "'16.13"
$ws.Range("E12").Value = "  +2.17%  "
$ws.Range("D13").Value = "'2.486.24"
$ws.Range("E13").Value = "  +3.00%  "
$ws.Range("D14").Value = "'22.14"
$ws.Range("E14").Value = "  -0.04%  "
$ws.Range("E15").Value = "  +0.91%  "
$ws.Range("E16").Value = "  +0.46%  "
$ws.Range("D17").Value = "'2.167.69"
$ws.Range("E17").Value = "  +3.02%  "
$ws.Range("D18").Value = "'39.468.17"
$ws.Range("E18").Value = "  +1.71%  "
$ws.Range("D19").Value = "'6.18"
$ws.Range("E19").Value = "  +1.05%  "
$ws.Range("D20").Value = "'71.91"
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("D21").Value = "'0.0₃0851"
$ws.Range("E21").Value = "  +1.14%  "
$ws.Range("D22").Value = "'229.71"
$ws.Range("E22").Value = "  +0.66%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").Value = "'2.34"
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").Value = "'2.36"
$ws.Range("E25").Value = "  +1.83%  "
$ws.Range("D26").Value = "'172.42"
$ws.Range("E26").Value = "  +0.40%  "
$ws.Range("D27").Value = "'9.53"
$ws.Range("E27").Value = "  -0.48%  "
$ws.Range("E28").Value = "  +1.73%  "
$ws.Range("D29").Value = "'19.89"
$ws.Range("E29").Value = "  +2.70%  "
$ws.Range("E30").Value = "  -0.36%  "
$ws.Range("D31").Value = "'2.62"
$ws.Range("E31").Value = "  +4.49%  "
$ws.Range("D33").Value = "'4.64"
$ws.Range("E33").Value = "  +1.93%  "
$ws.Range("E34").Value = "  +1.95%  "
$ws.Range("E35").Value = "  -0.70%  "
$ws.Range("D36").Value = "'0.0620"
$ws.Range("E36").Value = "  +0.32%  "
$ws.Range("E37").Value = "  +1.00%  "
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("D40").Value = "'103.21"
$ws.Range("E40").Value = "  +0.36%  "
$ws.Range("D41").Value = "'0.0230"
$ws.Range("E41").Value = "  +0.77%  "
$ws.Range("D42").Value = "'17.85"
$ws.Range("E42").Value = "  -1.45%  "
$ws.Range("D43").Value = "'1.524.07"
$ws.Range("E43").Value = "  -0.89%  "
$ws.Range("D44").Value = "'1.21"
$ws.Range("E44").Value = "  +3.99%  "
$ws.Range("E45").Value = "  +6.04%  "
$ws.Range("E46").Value = "  +1.89%  "
$ws.Range("E47").Value = "  +0.82%  "
$ws.Range("D48").Value = "'4.28"
$ws.Range("E48").Value = "  +3.93%  "
$ws.Range("E49").Value = "  -1.48%  "
$ws.Range("D50").Value = "'2.369.78"
$ws.Range("E50").Value = "  +3.02%  "
$ws.Range("D51").Value = "'2.95"
$ws.Range("E51").Value = "  -0.73%  "
